# Renumber the embedded Pearson/BTec logo pictures in the document's
# headers and footers:
#   - BTec_Logo-Orange pictures (in both headers):  image1.jpg -> image2.jpg
#   - PearsonLogo pictures      (in both footers):  image2.png -> image1.png
#
# The "name" seen here is the inline picture's internal docPr/name
# bookkeeping attribute (not visible document text), so it is reached
# through InlineShape.Name rather than Find/Replace.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineShape($range, $newName) {
    # Route the assignment through Selection - this is the reliable path
    # for both header and footer inline shapes in this object model.
    $range.Select() | Out-Null
    $shp = $word.Selection.InlineShapes.Item(1)
    $shp.Name = $newName
}

# Headers (BTec_Logo-Orange): image1.jpg -> image2.jpg
for ($i = 1; $i -le 2; $i++) {
    $hdr = $sec.Headers.Item($i)
    if ($hdr.Exists -and $hdr.Range.InlineShapes.Count -gt 0) {
        Rename-InlineShape $hdr.Range.InlineShapes.Item(1).Range "image2.jpg"
    }
}

# Footers (PearsonLogo): image2.png -> image1.png
for ($i = 1; $i -le 2; $i++) {
    $ftr = $sec.Footers.Item($i)
    if ($ftr.Exists -and $ftr.Range.InlineShapes.Count -gt 0) {
        Rename-InlineShape $ftr.Range.InlineShapes.Item(1).Range "image1.png"
    }
}
